# Applies the cryptos-list price/volume refresh described in the commit
# "Updated cryptos list on Sat Mar  2 04:40:16 UTC 2024 with GitHub Actions".
# Columns: D = Price, E = Volume(1h) (both stored as text in the sheet).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = '62.222.78'
$ws.Range("E2").Value = '  +1.28%  '
$ws.Range("D3").Value = '3.435.68'
$ws.Range("E3").Value = '  +1.66%  '
$ws.Range("E4").Value = '  -0.07%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '412.23'
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = '  +1.26%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.92'
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = '  -4.55%  '
$ws.Range("E7").Value = '  +6.01%  '
$ws.Range("E8").Value = '  -0.12%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.761'
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = '  +13.22%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.143'
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = '  +18.50%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '43.60'
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = '  +1.03%  '
$ws.Range("E12").Value = '  -0.23%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '8.99'
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = '  +7.25%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '20.74'
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = '  +5.29%  '
$ws.Range("E15").Value = '  +55.91%  '
$ws.Range("D16").Value = '3.449.34'
$ws.Range("E16").Value = '  +2.07%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '12.54'
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = '  +14.06%  '
$ws.Range("E18").Value = '  +4.95%  '
$ws.Range("D19").Value = '62.189.29'
$ws.Range("E19").Value = '  +1.30%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '405.26'
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = '  +29.23%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '90.20'
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = '  +8.29%  '
$ws.Range("E22").Value = '  -0.25%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '13.41'
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = '  +4.05%  '
$ws.Range("E24").Value = '  +3.25%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '33.79'
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = '  +14.62%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '8.74'
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = '  +4.88%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '4.81'
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = '  +0.44%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '7.70'
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = '  +0.00%  '
$ws.Range("E29").Value = '  +10.25%  '
$ws.Range("E30").Value = '  -0.10%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '44.07'
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = '  +6.91%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '11.96'
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = '  +5.56%  '
$ws.Range("E33").Value = '  +0.09%  '
$ws.Range("E34").Value = '  -0.07%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '0.0504'
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = '  +5.01%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '52.54'
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = '  +0.67%  '
$ws.Range("E37").Value = '  +0.07%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '3.42'
$ws.Range("D38").Style = "Normal"
$ws.Range("E38").Value = '  -0.54%  '
$ws.Range("E39").Value = '  -0.40%  '
$ws.Range("E40").Value = '  +6.99%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '0.315'
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = '  +5.17%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '140.86'
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = '  +2.17%  '
$ws.Range("E43").Value = '  +0.50%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.08'
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = '  -0.86%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '2.42'
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = '  +8.55%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '16.86'
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = '  +1.13%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '22.43'
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = '  +5.02%  '
$ws.Range("D48").Value = '2.126.52'
$ws.Range("E48").Value = '  -0.19%  '
$ws.Range("E49").Value = '  +1.75%  '
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = '0.0372'
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = '  +7.14%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '0.124'
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = '  +12.13%  '
